$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BIO - Landfills")

# --- Row 10: South Cardup Landfill Gas Project -----------------------------
$ws.Range("A10").Value = "South Cardup Landfill Gas Project"
$ws.Range("B10").Value = "WA"
$ws.Range("C10").Value = "LMS Energy"
$ws.Range("D10").Value = "site receiving between 18,000 and 25,000 tonnes of non-putrescible and putrescible waste per month"
$ws.Range("E10").Formula = "=21.5*12*1000"
$ws.Range("F10").Value = "MSW incineration"
$ws.Range("G10").Value = 0.9
$ws.Range("H10").Formula = "=E10*G10"
$ws.Range("I10").Value = "Mixed"
$ws.Range("J10").Value = "https://selectcivil.com.au/projects/south-cardup-landfill-sita-australia/"
$ws.Range("K10").Value = 116.01916872307901
$ws.Range("L10").Value = -32.266404199858698

# --- Row 11: Tamala Park LFG Power Station ----------------------------------
# (Values are written in the same order the shared-string table picks up new
# strings in the reference workbook: name, state, company, URL, then the
# capacity description.)
$ws.Range("A11").Value = "Tamala Park LFG Power Station"
$ws.Range("B11").Value = "WA"
$ws.Range("C11").Value = "EDL"
$ws.Range("F11").Value = "Landfill (LFG engines/flaring)"
$ws.Range("G11").Value = 0.8
$ws.Range("H11").Value = 235000
$ws.Range("I11").Value = "Yes (mostly)"
$ws.Range("J11").Value = "https://www.mrc.wa.gov.au/tamala-park/other-services/landfill-gas.aspx#:~:text=Tamala%20Park%20is%20home%20to,found%20on%20the%20EDL%20website."
$ws.Range("D11").Value = "capacity of six megawatts"
$ws.Range("K11").Value = 115.726560963035
$ws.Range("L11").Value = -31.707826453571599

# --- Hyperlinks (added before the format copy below so the final paste of
#     the row-9 template format is what determines the cells' visible style) -
$ws.Hyperlinks.Add($ws.Range("J10"), "https://selectcivil.com.au/projects/south-cardup-landfill-sita-australia/")
$ws.Hyperlinks.Add($ws.Range("J11"), "https://www.mrc.wa.gov.au/tamala-park/other-services/landfill-gas.aspx", ":~:text=Tamala%20Park%20is%20home%20to,found%20on%20the%20EDL%20website.")

# --- Copy the formatting from the template row (row 9) onto the new rows ---
$ws.Range("A9:L9").Copy()
$ws.Range("A10:L11").PasteSpecial(-4122)
# Row 11 column E has no value -- its format should match the other plain
# numeric cells (style used on row 2 etc.) rather than the row-9 template.
$ws.Range("E2").Copy()
$ws.Range("E11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Tab / selection state ----------------------------------------------------
$ws.Activate()
$ws.Range("A11").Select()

Write-Host "done"
